# Story sheet: rename header cells + convert embedded "\n" markers to "&" line-break markers
$wb = $excel.ActiveWorkbook
$story = $wb.Worksheets.Item("Story")
$message = $wb.Worksheets.Item("Message")

$story.Range("B1").Value = "Text"
$story.Range("C1").Value = "Delete"
$story.Range("B5").Value = "&&전 세계에 균열이 발생했다."
$story.Range("B6").Value = "차원 간의 균열이 열리며&튀어나온 온갖 마물들은&그야말로 자연재해였다."
$story.Range("B7").Value = "하지만,&인류는 멸망하지 않았다.&&균열에서 새어 나오는 마나를 받아들여 특별한 힘을 지니게 된 사람들 덕분이었다."
$story.Range("B8").Value = "내가 바로 그 힘을 가진&'헌터'다."
$story.Range("B9").Value = "나는 협회에 소속된 헌터다.&그것도....&헌터에 관련된 계약에 허점이 많은 시절에 계약하여 노예와 다름 없는 계약..."
$story.Range("B14").Value = "(사이렌 소리)하… 또 균열이 발생했다.&뭔 놈의 마물이 이렇게 하루가 멀다 하고 매일 나오는지 지겹다 지겨워…"
$story.Range("B15").Value = "(전화벨소리)왜?&(중얼거리며) 아니...학교에서 필요한 준비물이 있는데..&귀찮게 연락하지 말고, 카드로 사&(뚝 끊음)&"

# View-state: make Story the active/selected sheet with B18 selected,
# and leave Message no longer the tab-selected sheet (selection stays at A2).
$message.Range("A2").Select()
$story.Activate()
$story.Range("B18").Select()
